$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt9a"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3725013333333333
$ws.Range("H2").Value = 1.117504
$ws.Range("I2").Value = 0.05990362118104055
$ws.Range("J2").Value = 0.06743099450495174
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.524618
$ws.Range("N2").Value = 58.573854
$ws.Range("O2").Value = 0.4154885426712971
$ws.Range("P2").Value = 0.4539723485554654
$ws.Range("Q2").Value = 7.272946237824
$ws.Range("R2").Value = 65.45651614041599
$ws.Range("S2").Value = 0.02488926826524399
$ws.Range("T2").Value = 0.03061180694084363

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt9a"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3725013333333333
$ws.Range("H3").Value = 1.117504
$ws.Range("I3").Value = 0.05990362118104055
$ws.Range("J3").Value = 0.06743099450495174
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.24435933333334
$ws.Range("N3").Value = 45.73307800000001
$ws.Range("O3").Value = 0.324403614112412
$ws.Range("P3").Value = 0.3544508583357054
$ws.Range("Q3").Value = 5.678544177479111
$ws.Range("R3").Value = 51.106897597312
$ws.Range("S3").Value = 0.01943295120955039
$ws.Range("T3").Value = 0.02390097388071038

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt9a"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3725013333333333
$ws.Range("H4").Value = 1.117504
$ws.Range("I4").Value = 0.05990362118104055
$ws.Range("J4").Value = 0.06743099450495174
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.9507005
$ws.Range("N4").Value = 23.901401
$ws.Range("O4").Value = 0.2543137660693869
$ws.Range("P4").Value = 0.1852460510065796
$ws.Range("Q4").Value = 4.451651870517333
$ws.Range("R4").Value = 26.709911223104
$ws.Range("S4").Value = 0.01523431550374432
$ws.Range("T4").Value = 0.01249132544748868

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt9a"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3725013333333333
$ws.Range("H5").Value = 1.117504
$ws.Range("I5").Value = 0.05990362118104055
$ws.Range("J5").Value = 0.06743099450495174
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.272275
$ws.Range("N5").Value = 0.816825
$ws.Range("O5").Value = 0.005794077146903843
$ws.Range("P5").Value = 0.006330742102249548
$ws.Range("Q5").Value = 0.1014228005333333
$ws.Range("R5").Value = 0.9128052047999999
$ws.Range("S5").Value = 0.0003470862025018521
$ws.Range("T5").Value = 0.0004268882359090559

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt9a"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.763367333333334
$ws.Range("H6").Value = 11.290102
$ws.Range("I6").Value = 0.6052040917109096
$ws.Range("J6").Value = 0.6812528688240443
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.524618
$ws.Range("N6").Value = 58.573854
$ws.Range("O6").Value = 0.4154885426712971
$ws.Range("P6").Value = 0.4539723485554654
$ws.Range("Q6").Value = 73.478309577012
$ws.Range("R6").Value = 661.304786193108
$ws.Range("S6").Value = 0.2514553660836719
$ws.Range("T6").Value = 0.3092699648201998

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt9a"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.763367333333334
$ws.Range("H7").Value = 11.290102
$ws.Range("I7").Value = 0.6052040917109096
$ws.Range("J7").Value = 0.6812528688240443
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.24435933333334
$ws.Range("N7").Value = 45.73307800000001
$ws.Range("O7").Value = 0.324403614112412
$ws.Range("P7").Value = 0.3544508583357054
$ws.Range("Q7").Value = 57.37012393266179
$ws.Range("R7").Value = 516.3311153939561
$ws.Range("S7").Value = 0.1963303946266387
$ws.Range("T7").Value = 0.2414706640983442

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt9a"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.763367333333334
$ws.Range("H8").Value = 11.290102
$ws.Range("I8").Value = 0.6052040917109096
$ws.Range("J8").Value = 0.6812528688240443
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 11.9507005
$ws.Range("N8").Value = 23.901401
$ws.Range("O8").Value = 0.2543137660693869
$ws.Range("P8").Value = 0.1852460510065796
$ws.Range("Q8").Value = 44.97487587215033
$ws.Range("R8").Value = 269.849255232902
$ws.Range("S8").Value = 0.1539117318036041
$ws.Range("T8").Value = 0.1261994036865576

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt9a"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.763367333333334
$ws.Range("H9").Value = 11.290102
$ws.Range("I9").Value = 0.6052040917109096
$ws.Range("J9").Value = 0.6812528688240443
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.272275
$ws.Range("N9").Value = 0.816825
$ws.Range("O9").Value = 0.005794077146903843
$ws.Range("P9").Value = 0.006330742102249548
$ws.Range("Q9").Value = 1.024670840683333
$ws.Range("R9").Value = 9.222037566150002
$ws.Range("S9").Value = 0.003506599196994879
$ws.Range("T9").Value = 0.004312836218942666

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Wnt9a"
$ws.Range("C10").Value = "Fzd4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.0824755
$ws.Range("H10").Value = 4.164951
$ws.Range("I10").Value = 0.3348922871080498
$ws.Range("J10").Value = 0.2513161366710037
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 19.524618
$ws.Range("N10").Value = 58.573854
$ws.Range("O10").Value = 0.4154885426712971
$ws.Range("P10").Value = 0.4539723485554654
$ws.Range("Q10").Value = 40.659538631859
$ws.Range("R10").Value = 243.957231791154
$ws.Range("S10").Value = 0.1391439083223812
$ws.Range("T10").Value = 0.1140905767944219

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Wnt9a"
$ws.Range("C11").Value = "Fzd4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.0824755
$ws.Range("H11").Value = 4.164951
$ws.Range("I11").Value = 0.3348922871080498
$ws.Range("J11").Value = 0.2513161366710037
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 15.24435933333334
$ws.Range("N11").Value = 45.73307800000001
$ws.Range("O11").Value = 0.324403614112412
$ws.Range("P11").Value = 0.3544508583357054
$ws.Range("Q11").Value = 31.74600482486301
$ws.Range("R11").Value = 190.476028949178
$ws.Range("S11").Value = 0.1086402682762229
$ws.Range("T11").Value = 0.08907922035665071

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Wnt9a"
$ws.Range("C12").Value = "Fzd4"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.0824755
$ws.Range("H12").Value = 4.164951
$ws.Range("I12").Value = 0.3348922871080498
$ws.Range("J12").Value = 0.2513161366710037
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 11.9507005
$ws.Range("N12").Value = 23.901401
$ws.Range("O12").Value = 0.2543137660693869
$ws.Range("P12").Value = 0.1852460510065796
$ws.Range("Q12").Value = 24.88704099908775
$ws.Range("R12").Value = 99.54816399635101
$ws.Range("S12").Value = 0.08516771876203855
$ws.Range("T12").Value = 0.0465553218725333

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Wnt9a"
$ws.Range("C13").Value = "Fzd4"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.0824755
$ws.Range("H13").Value = 4.164951
$ws.Range("I13").Value = 0.3348922871080498
$ws.Range("J13").Value = 0.2513161366710037
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.272275
$ws.Range("N13").Value = 0.816825
$ws.Range("O13").Value = 0.005794077146903843
$ws.Range("P13").Value = 0.006330742102249548
$ws.Range("Q13").Value = 0.5670060167625001
$ws.Range("R13").Value = 3.402036100575
$ws.Range("S13").Value = 0.001940391747407112
$ws.Range("T13").Value = 0.001591017647397825

Write-Output "done"